$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.906

$ws.Range("A9").Value = -21.882
$ws.Range("C9").Value = -12.222
$ws.Range("D9").Value = -7.866

$ws.Range("A18").Value = -21.898

$ws.Range("A20").Value = -20.563

$ws.Range("C23").Value = -12.676

$ws.Range("C24").Value = -12.189

$ws.Range("C26").Value = -12.401

$ws.Range("A27").Value = -21.888

$ws.Range("D32").Value = -7.382000000000001

$ws.Range("C34").Value = -11.853

$ws.Range("C35").Value = -12.415

$ws.Range("D38").Value = -7.905999999999999

$ws.Range("D45").Value = -7.458999999999999

$ws.Range("C48").Value = -11.529

$ws.Range("D51").Value = -8.059999999999999

$ws.Range("C52").Value = -11.749

$ws.Range("D57").Value = -8.145999999999999

$ws.Range("D64").Value = -7.711

$ws.Range("C66").Value = -11.574

$ws.Range("C67").Value = -10.977

$ws.Range("A69").Value = -21.593

$ws.Range("A76").Value = -20.287

$ws.Range("C80").Value = -12.522

$ws.Range("A82").Value = -22.004

$ws.Range("D93").Value = -7.052

$ws.Range("C99").Value = -11.695
